# Apply odds updates to Sheet1, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("S2").Value = 1.53
$ws.Range("T2").Value = 2.38

# Row 3
$ws.Range("S3").Value = 1.53
$ws.Range("T3").Value = 2.38

# Row 4
$ws.Range("S4").Value = 1.57

# Row 5
$ws.Range("I5").Value = 3.6
$ws.Range("K5").Value = 2.05
$ws.Range("L5").Value = 4.33
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 9
$ws.Range("O5").Value = 1.4
$ws.Range("P5").Value = 2.75
$ws.Range("Q5").Value = 2.2
$ws.Range("R5").Value = 1.65
$ws.Range("S5").Value = 1.5
$ws.Range("U5").Value = 2
$ws.Range("V5").Value = 1.73
$ws.Range("W5").Value = 6.5
$ws.Range("X5").Value = 9
$ws.Range("AC5").Value = 8
$ws.Range("AG5").Value = 9
$ws.Range("AP5").Value = 23
$ws.Range("AU5").Value = 8.5
$ws.Range("AX5").Value = 21
$ws.Range("BA5").Value = 101

# Row 9
$ws.Range("J9").Value = 3.75
$ws.Range("S9").Value = 1.4
$ws.Range("T9").Value = 2.75
$ws.Range("U9").Value = 1.75
$ws.Range("V9").Value = 2
$ws.Range("AB9").Value = 34
$ws.Range("AE9").Value = 15
$ws.Range("AF9").Value = 51
$ws.Range("AG9").Value = 8
$ws.Range("AP9").Value = 26
$ws.Range("AT9").Value = 2.75

# Row 10
$ws.Range("G10").Value = 2.25
$ws.Range("I10").Value = 3
$ws.Range("K10").Value = 2.25
$ws.Range("U10").Value = 1.62
$ws.Range("V10").Value = 2.2
$ws.Range("W10").Value = 9.5
$ws.Range("X10").Value = 12
$ws.Range("AE10").Value = 13
$ws.Range("AL10").Value = 29
$ws.Range("AP10").Value = 21
$ws.Range("BA10").Value = 67
